$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source inlineStr cells)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '60.403.94'
$ws.Range("E2").Value = '  +2.47%  '
$ws.Range("D3").Value = '2.695.57'
$ws.Range("E3").Value = '  +1.40%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '523.69'
$ws.Range("E5").Value = '  +1.94%  '
$ws.Range("D6").Value = '147.05'
$ws.Range("E6").Value = '  +2.17%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +2.18%  '
$ws.Range("D9").Value = '2.718.36'
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("D10").Value = '6.49'
$ws.Range("E10").Value = '  +3.43%  '
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '0.341'
$ws.Range("E12").Value = '  +1.72%  '
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").Value = '3.176.82'
$ws.Range("E14").Value = '  +1.57%  '
$ws.Range("D15").Value = '60.465.41'
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.836.62'
$ws.Range("E16").Value = '  +5.44%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '21.43'
$ws.Range("E17").Value = '  +1.89%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.0000139'
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").Value = '352.94'
$ws.Range("E19").Value = '  +0.69%  '
$ws.Range("D20").Value = '4.55'
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("D21").Value = '10.62'
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("E22").Value = '  +3.80%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '62.98'
$ws.Range("E24").Value = '  +3.24%  '
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("E26").Value = '  +5.50%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").Value = '0.0₃0822'
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").Value = '7.36'
$ws.Range("E29").Value = '  +1.27%  '
$ws.Range("D30").Value = '6.90'
$ws.Range("E30").Value = '  +7.08%  '
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  +1.60%  '
$ws.Range("D33").Value = '19.19'
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("D34").Value = '147.83'
$ws.Range("E34").Value = '  -1.33%  '
$ws.Range("D35").Value = '4.31'
$ws.Range("E35").Value = '  +7.73%  '
$ws.Range("E36").Value = '  +9.52%  '
$ws.Range("D37").Value = '0.955'
$ws.Range("E37").Value = '  -5.84%  '
$ws.Range("D38").Value = '1.54'
$ws.Range("E38").Value = '  +9.89%  '
$ws.Range("D39").Value = '0.879'
$ws.Range("E39").Value = '  +3.92%  '
$ws.Range("D40").Value = '36.93'
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").Value = '285.56'
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").Value = '20.22'
$ws.Range("E43").Value = '  +2.31%  '
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D45").Value = '0.0991'
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = '0.996'
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.143.76'
$ws.Range("E47").Value = '  +6.94%  '
$ws.Range("E48").Value = '  +4.69%  '
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").Value = '0.0235'
$ws.Range("E50").Value = '  +2.22%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = '10.46'
$ws.Range("E51").Value = '  +1.88%  '
